# SOHP Excel template update
# Adds "Interviewee Date" / "Interviewer Date" name-part columns, and
# consolidates the separate "Interview Number" note into the identifier
# element (with a displayLabel attribute), and appends an LCGFT genre
# tag to the typeOfResource cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 columns before F for the Interviewee "date" namePart ---
# (old F:W shift right by 3 -> old F now lives at I)
$ws.Columns("F:H").Insert()

# --- 2. Insert 3 columns before the (now-shifted) old I for the
#        Interviewer "date" namePart. Old I is currently at column L. ---
$ws.Columns("L:N").Insert()

# --- 3. Remove the old standalone "Interview Number" note block, which
#        (after the two inserts above) now sits at Y:AA. Its content is
#        being folded into the identifier element instead. ---
$ws.Columns("Y:AA").Delete()

# --- 4. Fill in the new Interviewee-date cells (F1:H1) ---
$ws.Range("F1").Value = "</mods:namePart>"
$ws.Range("G1").Value = "<mods:namePart type=""date"">"
$ws.Range("H1").Value = "Interviewee Date"
$ws.Range("F1").Font.Bold = $false
$ws.Range("G1").Font.Bold = $false
$ws.Range("H1").Font.Bold = $true

# --- 5. Fill in the new Interviewer-date cells (L1:N1) ---
$ws.Range("L1").Value = "</mods:namePart>"
$ws.Range("M1").Value = "<mods:namePart type=""date"">"
$ws.Range("N1").Value = "Interviewer Date"
$ws.Range("L1").Font.Bold = $false
$ws.Range("M1").Font.Bold = $false
$ws.Range("N1").Font.Bold = $true

# --- 6. Update the identifier open tag + label (now at V1:W1) so the
#        "Interview Number" displayLabel lives on the identifier itself ---
$ws.Range("V1").Value = "<mods:identifier displayLabel=""Interview Number"" type=""local"">"
$ws.Range("W1").Value = "Interview Number"

# --- 7. Append the LCGFT "Oral histories" genre to typeOfResource (now Y1) ---
$ws.Range("Y1").Value = "<mods:typeOfResource>sound recording-nonmusical</mods:typeOfResource><mods:genre authority=""lcgft"">Oral histories</mods:genre>"

# --- 8. Scroll/selection so column E is the leftmost visible column and
#        Z1 (the final populated cell) is selected, matching the saved view ---
$ws.Range("Z1").Select()
$excel.ActiveWindow.ScrollColumn = 5
